# Write sum of donations made to the nominator per committee-member to a
# new worksheet ("sum_donations_to_nominator").

$wb = $excel.ActiveWorkbook
$fullNameMatch = $wb.Worksheets.Item("full_name_match")

# The "full_name_match" sheet's header row mislabeled columns D/E - fix
# the labels (D = Campaign, E = Nominated by) to match their actual data.
$fullNameMatch.Range("D1").Value = "Campaign"
$fullNameMatch.Range("E1").Value = "Nominated by"

# Pull every row (including the now-corrected header) in one shot.
$lastRow = $fullNameMatch.Cells.Item($fullNameMatch.Rows.Count, 1).End(-4162).Row
$vals = $fullNameMatch.Range("A1:E" + $lastRow).Value()
$rowCount = $vals.GetLength(0)

# Walk the donation rows and keep only the ones where the donation's
# "Nominated by" full name actually resolves to the "Campaign" last name
# recorded alongside it (i.e. a genuine nominator match). Sum the donation
# Amount per committee member, per contiguous block sharing the same
# (Committee Member, Campaign) pair - mirrors how the source rows are laid
# out (a committee member's matching donations are grouped together).
$names = @()
$sums = @()
$camps = @()
$curKey = $null
$curIdx = -1

for ($r = 2; $r -le $rowCount; $r++) {
    $committee = $vals[$r, 1]
    $amount = [double]$vals[$r, 3]
    $nominatedBy = $vals[$r, 4]
    $campaign = $vals[$r, 5]

    $parts = $nominatedBy.Split(" ")
    $lastName = $parts[$parts.Length - 1].ToUpper()

    if ($lastName -eq $campaign) {
        $key = $committee + "|" + $campaign
        if ($key -eq $curKey) {
            $sums[$curIdx] = $sums[$curIdx] + $amount
        } else {
            $names += $committee
            $sums += $amount
            $camps += $campaign
            $curIdx = $curIdx + 1
            $curKey = $key
        }
    }
}

$outRows = $names.Length
$outDim0 = $outRows + 1

# Build the new worksheet, positioned after "full_name_match".
$newSheet = $wb.Worksheets.Add($null, $fullNameMatch)
$newSheet.Name = "sum_donations_to_nominator"

$out = New-Object 'object[,]' $outDim0, 3
$out[0, 0] = "Committee Member"
$out[0, 1] = "Amount"
$out[0, 2] = "Campaign"

for ($i = 0; $i -lt $outRows; $i++) {
    $out[$i + 1, 0] = $names[$i]
    $out[$i + 1, 1] = $sums[$i]
    $out[$i + 1, 2] = $camps[$i]
}

$endCell = $newSheet.Cells.Item($outRows + 1, 3)
$newSheet.Range($newSheet.Cells.Item(1, 1), $endCell).Value = $out

$wb.Worksheets.Item("last_name_match").Select()
